# Apply the "backup" column + corrected detect_structure values + six new
# monthly rows to the NCC.NS stock-history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column R ("backup"): header styled like the other headers, and
#    0 for every existing data row (2-263).
# ---------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value2 = "backup"

$ws.Range("R2:R263").Value2 = 0

# ---------------------------------------------------------------------
# 2. Three stale "detect_structure" flags got reset to 0 once the backup
#    column took over tracking those events.
# ---------------------------------------------------------------------
$ws.Range("Q26").Value2 = 0
$ws.Range("Q49").Value2 = 0
$ws.Range("Q51").Value2 = 0

# ---------------------------------------------------------------------
# 3. Six new monthly rows (264-269), continuing the series through
#    Dec 2024. Column A keeps the same datetime number format as the
#    rest of the column; F (Adj Close) and R (backup) are left blank,
#    matching source rows that had no adjusted-close data yet.
# ---------------------------------------------------------------------
$ws.Range("A264:A269").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(264, 1).Value2 = 45474
$ws.Cells.Item(264, 2).Value2 = 314.3233873342549
$ws.Cells.Item(264, 3).Value2 = 361.9932849394499
$ws.Cells.Item(264, 4).Value2 = 293.4678071319821
$ws.Cells.Item(264, 5).Value2 = 355.4883422851562
$ws.Cells.Item(264, 7).Value2 = 197706340
$ws.Cells.Item(264, 8).Value2 = 2024
$ws.Cells.Item(264, 9).Value2 = 7
$ws.Cells.Item(264, 10).Value2 = 1
$ws.Cells.Item(264, 11).Value2 = 0
$ws.Cells.Item(264, 12).Value2 = 0
$ws.Cells.Item(264, 13).Value2 = 0
$ws.Cells.Item(264, 14).Value2 = 27
$ws.Cells.Item(264, 15).Value2 = 1
$ws.Cells.Item(264, 16).Value2 = 0
$ws.Cells.Item(264, 17).Value2 = 0

$ws.Cells.Item(265, 1).Value2 = 45505
$ws.Cells.Item(265, 2).Value2 = 357.1269849407552
$ws.Cells.Item(265, 3).Value2 = 360.1063535076007
$ws.Cells.Item(265, 4).Value2 = 306.3784009572792
$ws.Cells.Item(265, 5).Value2 = 316.8558349609375
$ws.Cells.Item(265, 7).Value2 = 111455028
$ws.Cells.Item(265, 8).Value2 = 2024
$ws.Cells.Item(265, 9).Value2 = 8
$ws.Cells.Item(265, 10).Value2 = 1
$ws.Cells.Item(265, 11).Value2 = 0
$ws.Cells.Item(265, 12).Value2 = 0
$ws.Cells.Item(265, 13).Value2 = 0
$ws.Cells.Item(265, 14).Value2 = 31
$ws.Cells.Item(265, 15).Value2 = 0
$ws.Cells.Item(265, 16).Value2 = 0
$ws.Cells.Item(265, 17).Value2 = 0

$ws.Cells.Item(266, 1).Value2 = 45536
$ws.Cells.Item(266, 2).Value2 = 322.5
$ws.Cells.Item(266, 3).Value2 = 330.7999877929688
$ws.Cells.Item(266, 4).Value2 = 293.0499877929688
$ws.Cells.Item(266, 5).Value2 = 301.75
$ws.Cells.Item(266, 7).Value2 = 61112586
$ws.Cells.Item(266, 8).Value2 = 2024
$ws.Cells.Item(266, 9).Value2 = 9
$ws.Cells.Item(266, 10).Value2 = 1
$ws.Cells.Item(266, 11).Value2 = 0
$ws.Cells.Item(266, 12).Value2 = 0
$ws.Cells.Item(266, 13).Value2 = 0
$ws.Cells.Item(266, 14).Value2 = 35
$ws.Cells.Item(266, 15).Value2 = 0
$ws.Cells.Item(266, 16).Value2 = 0
$ws.Cells.Item(266, 17).Value2 = 0

$ws.Cells.Item(267, 1).Value2 = 45566
$ws.Cells.Item(267, 2).Value2 = 301.75
$ws.Cells.Item(267, 3).Value2 = 312.8999938964844
$ws.Cells.Item(267, 4).Value2 = 270.6499938964844
$ws.Cells.Item(267, 5).Value2 = 298.7999877929688
$ws.Cells.Item(267, 7).Value2 = 63126083
$ws.Cells.Item(267, 8).Value2 = 2024
$ws.Cells.Item(267, 9).Value2 = 10
$ws.Cells.Item(267, 10).Value2 = 1
$ws.Cells.Item(267, 11).Value2 = 0
$ws.Cells.Item(267, 12).Value2 = 0
$ws.Cells.Item(267, 13).Value2 = 0
$ws.Cells.Item(267, 14).Value2 = 40
$ws.Cells.Item(267, 15).Value2 = 0
$ws.Cells.Item(267, 16).Value2 = 0
$ws.Cells.Item(267, 17).Value2 = 0

$ws.Cells.Item(268, 1).Value2 = 45597
$ws.Cells.Item(268, 2).Value2 = 317.9500122070312
$ws.Cells.Item(268, 3).Value2 = 321.8999938964844
$ws.Cells.Item(268, 4).Value2 = 269.25
$ws.Cells.Item(268, 5).Value2 = 309.7000122070312
$ws.Cells.Item(268, 7).Value2 = 62109874
$ws.Cells.Item(268, 8).Value2 = 2024
$ws.Cells.Item(268, 9).Value2 = 11
$ws.Cells.Item(268, 10).Value2 = 1
$ws.Cells.Item(268, 11).Value2 = 0
$ws.Cells.Item(268, 12).Value2 = 0
$ws.Cells.Item(268, 13).Value2 = 0
$ws.Cells.Item(268, 14).Value2 = 44
$ws.Cells.Item(268, 15).Value2 = 0
$ws.Cells.Item(268, 16).Value2 = 0
$ws.Cells.Item(268, 17).Value2 = 2

$ws.Cells.Item(269, 1).Value2 = 45627
$ws.Cells.Item(269, 2).Value2 = 311.8999938964844
$ws.Cells.Item(269, 3).Value2 = 326.4500122070312
$ws.Cells.Item(269, 4).Value2 = 268.5
$ws.Cells.Item(269, 5).Value2 = 270.2999877929688
$ws.Cells.Item(269, 7).Value2 = 56939033
$ws.Cells.Item(269, 8).Value2 = 2024
$ws.Cells.Item(269, 9).Value2 = 12
$ws.Cells.Item(269, 10).Value2 = 1
$ws.Cells.Item(269, 11).Value2 = 0
$ws.Cells.Item(269, 12).Value2 = 0
$ws.Cells.Item(269, 13).Value2 = 0
$ws.Cells.Item(269, 14).Value2 = 48
$ws.Cells.Item(269, 15).Value2 = 0
$ws.Cells.Item(269, 16).Value2 = 0
$ws.Cells.Item(269, 17).Value2 = 0

Write-Host "Done"
